$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain decimal numbers but must remain text
# (matching the source data, which stores prices as inline strings).
# Mark them as Text-formatted before assigning so Excel does not
# silently convert them to numeric values.
$textCells = @("D5","D7","D9","D10","D11","D12","D13","D15","D16","D20","D21","D23","D24","D26","D29","D30","D34","D35","D39","D42","D43","D44","D47","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin values scraped by the GitHub Actions job
$ws.Range('D2').Value = '42.332.77'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '2.244.85'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '246.05'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').Value = '75.65'
$ws.Range('E7').Value = '  -0.52%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.619'
$ws.Range('E9').Value = '  -2.10%  '
$ws.Range('D10').Value = '43.94'
$ws.Range('E10').Value = '  +8.81%  '
$ws.Range('D11').Value = '0.0947'
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').Value = '7.23'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').Value = '0.102'
$ws.Range('E13').Value = '  -2.00%  '
$ws.Range('D14').Value = '2.579.18'
$ws.Range('E14').Value = '  +0.02%  '
$ws.Range('D15').Value = '14.53'
$ws.Range('E15').Value = '  -2.41%  '
$ws.Range('D16').Value = '0.856'
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('D17').Value = '2.255.09'
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('D18').Value = '42.183.73'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('E19').Value = '  +4.43%  '
$ws.Range('D20').Value = '6.17'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').Value = '72.12'
$ws.Range('E21').Value = '  +0.82%  '
$ws.Range('E22').Value = '  +0.76%  '
$ws.Range('D23').Value = '230.05'
$ws.Range('E23').Value = '  -0.73%  '
$ws.Range('D24').Value = '9.22'
$ws.Range('E24').Value = '  +30.52%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').Value = '11.47'
$ws.Range('E26').Value = '  +2.80%  '
$ws.Range('E27').Value = '  -2.85%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = '2.19'
$ws.Range('E29').Value = '  +1.54%  '
$ws.Range('D30').Value = '167.93'
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('E32').Value = '  -3.36%  '
$ws.Range('E33').Value = '  +0.81%  '
$ws.Range('D34').Value = '30.83'
$ws.Range('E34').Value = '  -5.54%  '
$ws.Range('D35').Value = '5.33'
$ws.Range('E35').Value = '  +11.10%  '
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('E37').Value = '  -0.22%  '
$ws.Range('E38').Value = '  +5.76%  '
$ws.Range('D39').Value = '14.00'
$ws.Range('E39').Value = '  +4.45%  '
$ws.Range('E40').Value = '  -1.22%  '
$ws.Range('E41').Value = '  -2.38%  '
$ws.Range('D42').Value = '64.28'
$ws.Range('E42').Value = '  +6.88%  '
$ws.Range('D43').Value = '0.201'
$ws.Range('E43').Value = '  -0.79%  '
$ws.Range('D44').Value = '107.78'
$ws.Range('E44').Value = '  -8.34%  '
$ws.Range('E45').Value = '  +1.03%  '
$ws.Range('E46').Value = '  +1.83%  '
$ws.Range('D47').Value = '0.996'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Value = '1.13'
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').Value = '1.19'
$ws.Range('E49').Value = '  +1.27%  '
$ws.Range('E50').Value = '  +4.36%  '
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').Value = '2.71'
$ws.Range('E51').Value = '  +0.89%  '
